$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FN_PASS")

$ws.Range("A1").Value = "Karina"
$ws.Range("B1").Value = "adytx1yloe96so"
$ws.Range("A2").Value = "Melony"
$ws.Range("B2").Value = "ome71tlpz"
$ws.Range("A3").Value = "Dawn"
$ws.Range("B3").Value = "5qcrnbwc7"
